$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row at position 4 for the "family with disabilities
#    Persons" series (pushes the old row 4 -> 5 and old row 5 -> 6).
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# Borrow the plain numeric formatting (no border / no special alignment)
# that is already used by the un-bordered data cells so the new row matches
# the look of the other data rows.
$ws.Range("B5").Copy()
$ws.Range("B4:I4").PasteSpecial(-4122)
$ws.Range("B5:H5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Label + values for the new row 4 ("family with disabilities Persons")
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("B4").Value = 1279
$ws.Range("C4").Value = 1272
$ws.Range("D4").Value = 1289
$ws.Range("E4").Value = 1305
$ws.Range("F4").Value = 1322
$ws.Range("G4").Value = 1313
$ws.Range("H4").Value = 1263
$ws.Range("I4").Value = 1268

# Label + values for row 5 ("disabilities Persons" -- was the old row 4,
# "Number of disability persons")
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("B5").Value = 1426
$ws.Range("C5").Value = 1419
$ws.Range("D5").Value = 1432
$ws.Range("E5").Value = 1445
$ws.Range("F5").Value = 1460
$ws.Range("G5").Value = 1455
$ws.Range("H5").Value = 1401
$ws.Range("I5").Value = 1410

# Borders: row4 label gets a top border, row5 label gets a bottom border,
# and I5 keeps the bottom border that the former I4 ("Number of disability
# persons") cell had.
$ws.Range("A4").Borders.Item(9).LineStyle = 1
$ws.Range("A4").Borders.Item(9).Weight = 2
$ws.Range("A5").Borders.Item(10).LineStyle = 1
$ws.Range("A5").Borders.Item(10).Weight = 2
$ws.Range("I5").Borders.Item(10).LineStyle = 1
$ws.Range("I5").Borders.Item(10).Weight = 2
$ws.Range("I5").HorizontalAlignment = 1

$ws.Rows.Item(4).RowHeight = 24.75
$ws.Rows.Item(5).RowHeight = 21

# ---------------------------------------------------------------------------
# 2) Turn row 1 into the new merged title row (replaces the old title text).
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Merge()
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Chiatura Municipality"
$ws.Range("A1:I1").Font.Bold = $true
$ws.Range("A1:I1").Font.Name = "Arial"
$ws.Range("A1:I1").Font.Size = 11
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 51

# Column width for the label column.
$ws.Columns.Item(1).ColumnWidth = 19.98

$ws.Range("A1").Select()
